# Swap the "step" contents of TC2 and TC4 (TC3 stays the same).
# Before:
#   TC2 step2 (row 20): B="Chefe Clica para realizar a autorização de pagamento."
#                        D="SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento"
#   TC4 step2 (row 36): B="Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
#                        D="SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."
# After:
#   TC2 step2 (row 20): gets the old TC4 content
#   TC4 step2 (row 36): gets the old TC2 content
#   TC3 (rows 27-28): unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc2StepB = $ws.Range("B20").Value2
$tc2StepD = $ws.Range("D20").Value2

$tc4StepB = $ws.Range("B36").Value2
$tc4StepD = $ws.Range("D36").Value2

$ws.Range("B20").Value2 = $tc4StepB
$ws.Range("D20").Value2 = $tc4StepD

$ws.Range("B36").Value2 = $tc2StepB
$ws.Range("D36").Value2 = $tc2StepD
